# [5500] Timers and small fixes
#
# 1. HvN-Treinseries: sort the data range A2:F61 ascending by column D
#    (Consist Type) and leave the selection on B23.
# 2. Scenarios: mark row 2 (A2:C2) as "done" (same highlight style as the
#    other finished columns on that row) and fill in the Timer info for
#    row 3 (D3/G3/H3) by copying the equivalent cells from row 2.
# 3. Scenarios!F3: strike through the "- 1700/11700 Ut9" line, which is
#    now covered by the 1700/11700 timer fix noted above.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. HvN-Treinseries: re-sort by Consist Type (column D)
# ---------------------------------------------------------------------
$wsHvN = $wb.Worksheets.Item("HvN-Treinseries")
$wsHvN.Activate()

$sortRange = $wsHvN.Range("A1:F61")
$sortKey = $wsHvN.Range("D2:D61")
$sortRange.Sort($sortKey, 1)

$wsHvN.Range("B23").Select()

# ---------------------------------------------------------------------
# 2. Scenarios: highlight row 2 + fill in row 3 timer cells
# ---------------------------------------------------------------------
$wsScenarios = $wb.Worksheets.Item("Scenarios")
$wsScenarios.Activate()

# A2:C2 get the same "confirmed" highlight fill already used on D2/I2
$wsScenarios.Range("D2").Copy()
$wsScenarios.Range("A2:C2").PasteSpecial(-4122)  # xlPasteFormats

# D3 = "yes" (same as D2)
$wsScenarios.Range("D2").Copy()
$wsScenarios.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$wsScenarios.Range("D2").Copy()
$wsScenarios.Range("D3").PasteSpecial(-4163)  # xlPasteValues

# G3 = "- Ut" (same as G2)
$wsScenarios.Range("G2").Copy()
$wsScenarios.Range("G3").PasteSpecial(-4122)
$wsScenarios.Range("G2").Copy()
$wsScenarios.Range("G3").PasteSpecial(-4163)

# H3 = "yes" (same as D2)
$wsScenarios.Range("D2").Copy()
$wsScenarios.Range("H3").PasteSpecial(-4122)
$wsScenarios.Range("D2").Copy()
$wsScenarios.Range("H3").PasteSpecial(-4163)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Scenarios!F3: strike through the now-superseded "Ut9" timer line
# ---------------------------------------------------------------------
$f3 = $wsScenarios.Range("F3")
$chars = $f3.Characters(193, 16)
$chars.Font.Strikethrough = $true

$wsScenarios.Range("G3").Select()
